$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 180.08333
$ws.Range("I2").Value = 80.42856999999999
$ws.Range("J2").Value = 319.6
$ws.Range("K2").Value = 80.42856999999999
$ws.Range("L2").Value = 319.6
$ws.Range("M2").Value = 32.57143000000001
$ws.Range("N2").Value = -545.6
$ws.Range("H8").Value = 180
$ws.Range("I8").Value = 180
$ws.Range("K8").Value = 540
$ws.Range("M8").Value = -401
$ws.Range("H17").Value = 1953.7
$ws.Range("J17").Value = 1953.7
$ws.Range("L17").Value = 5861.1
$ws.Range("N17").Value = -6197.1
$ws.Range("H21").Value = 18000
$ws.Range("I21").Value = 18000
$ws.Range("K21").Value = 18000
$ws.Range("M21").Value = -17532
$ws.Range("H23").Value = 18000
$ws.Range("I23").Value = 18000
$ws.Range("K23").Value = 18000
$ws.Range("M23").Value = -17766
$ws.Range("H29").Value = 5
$ws.Range("I29").Value = 5
$ws.Range("K29").Value = 15
$ws.Range("M29").Value = 266
$ws.Range("H49").Value = 25
$ws.Range("I49").Value = 25
$ws.Range("K49").Value = 75
$ws.Range("M49").Value = 61
$ws.Range("H97").Value = 94999.5
$ws.Range("J97").Value = 94999.5
$ws.Range("L97").Value = 284998.5
$ws.Range("N97").Value = -285990.5
$ws.Range("H135").Value = 1081.5
$ws.Range("I135").Value = 1133.1578
$ws.Range("J135").Value = 100
$ws.Range("K135").Value = 10198.4202
$ws.Range("L135").Value = 900
$ws.Range("M135").Value = -7663.4202
$ws.Range("N135").Value = -5970

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11512.533
$ws.Range("I32").Value = 4672.606
$ws.Range("K32").Value = 4672.606
$ws.Range("M32").Value = -4385.606
$ws.Range("H42").Value = 88998
$ws.Range("J42").Value = 88998
$ws.Range("L42").Value = 88998
$ws.Range("N42").Value = -89970
$ws.Range("H61").Value = 64565.875
$ws.Range("I61").Value = 1829.5834
$ws.Range("K61").Value = 1829.5834
$ws.Range("M61").Value = -1617.5834
$ws.Range("H114").Value = 71700
$ws.Range("J114").Value = 71700
$ws.Range("L114").Value = 71700
$ws.Range("N114").Value = -80378
$ws.Range("H136").Value = 64565.875
$ws.Range("I136").Value = 1829.5834
$ws.Range("K136").Value = 5488.7502
$ws.Range("M136").Value = -2938.7502

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 30006.75
$ws.Range("J106").Value = 30006.75
$ws.Range("L106").Value = 30006.75
$ws.Range("N106").Value = -32530.75
$ws.Range("H107").Value = 3772.3076
$ws.Range("I107").Value = 2529.5
$ws.Range("K107").Value = 2529.5
$ws.Range("M107").Value = -609.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = $null
$ws.Range("N22").Value = 0
$ws.Range("H134").Value = 43790.457
$ws.Range("I134").Value = 2024.1
$ws.Range("K134").Value = 6072.299999999999
$ws.Range("M134").Value = -3537.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4703.3213
$ws.Range("J39").Value = 6006.6665
$ws.Range("L39").Value = 18019.9995
$ws.Range("N39").Value = -18607.9995
$ws.Range("H46").Value = 4999.5
$ws.Range("I46").Value = 9000
$ws.Range("J46").Value = 999
$ws.Range("K46").Value = 27000
$ws.Range("L46").Value = 2997
$ws.Range("M46").Value = -26909
$ws.Range("N46").Value = -3179
$ws.Range("H51").Value = 3105.375
$ws.Range("I51").Value = 1970.6
$ws.Range("J51").Value = 4996.6665
$ws.Range("K51").Value = 5911.799999999999
$ws.Range("L51").Value = 14989.9995
$ws.Range("M51").Value = -5451.799999999999
$ws.Range("N51").Value = -15909.9995
$ws.Range("H58").Value = 2413
$ws.Range("I58").Value = 619.5
$ws.Range("J58").Value = 6000
$ws.Range("K58").Value = 1858.5
$ws.Range("L58").Value = 18000
$ws.Range("M58").Value = -1730.5
$ws.Range("N58").Value = -18256
$ws.Range("H113").Value = 518
$ws.Range("I113").Value = 285.33334
$ws.Range("J113").Value = 667.5714
$ws.Range("K113").Value = 856.0000200000001
$ws.Range("L113").Value = 2002.7142
$ws.Range("M113").Value = 1313.99998
$ws.Range("N113").Value = -6342.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 40022500
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").Value = $null
$ws.Range("H31").Value = 9833.333000000001
$ws.Range("I31").Value = 750
$ws.Range("J31").Value = 28000
$ws.Range("K31").Value = 750
$ws.Range("L31").Value = 28000
$ws.Range("M31").Value = -458
$ws.Range("N31").Value = -28584
$ws.Range("H37").Value = 9833.333000000001
$ws.Range("I37").Value = 750
$ws.Range("J37").Value = 28000
$ws.Range("K37").Value = 750
$ws.Range("L37").Value = 28000
$ws.Range("M37").Value = -473
$ws.Range("N37").Value = -28554
$ws.Range("H40").Value = 6003.75
$ws.Range("I40").Value = 2015
$ws.Range("J40").Value = 7333.3335
$ws.Range("K40").Value = 2015
$ws.Range("L40").Value = 7333.3335
$ws.Range("M40").Value = -1864
$ws.Range("N40").Value = -7635.3335
$ws.Range("H80").Value = 2554.0667
$ws.Range("I80").Value = 2500.7144
$ws.Range("J80").Value = 2600.75
$ws.Range("K80").Value = 2500.7144
$ws.Range("L80").Value = 2600.75
$ws.Range("M80").Value = -1502.7144
$ws.Range("N80").Value = -4596.75
$ws.Range("H83").Value = 2554.0667
$ws.Range("I83").Value = 2500.7144
$ws.Range("J83").Value = 2600.75
$ws.Range("K83").Value = 12503.572
$ws.Range("L83").Value = 13003.75
$ws.Range("M83").Value = -7511.572
$ws.Range("N83").Value = -22987.75
$ws.Range("H86").Value = 59999.25
$ws.Range("J86").Value = 59999.25
$ws.Range("L86").Value = 59999.25
$ws.Range("N86").Value = -62371.25
$ws.Range("H89").Value = 59999.25
$ws.Range("J89").Value = 59999.25
$ws.Range("L89").Value = 179997.75
$ws.Range("N89").Value = -191853.75
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = $null
$ws.Range("N103").Value = 0
$ws.Range("H107").Value = 393.2857
$ws.Range("I107").Value = 459.5
$ws.Range("J107").Value = 343.625
$ws.Range("K107").Value = 459.5
$ws.Range("L107").Value = 343.625
$ws.Range("M107").Value = 1460.5
$ws.Range("N107").Value = -4183.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1438.6428
$ws.Range("I16").Value = 1286.4546
$ws.Range("K16").Value = 1286.4546
$ws.Range("M16").Value = -1116.4546
$ws.Range("H20").Value = 2516250
$ws.Range("J20").Value = 2516250
$ws.Range("L20").Value = 2516250
$ws.Range("N20").Value = -2516702
$ws.Range("H106").Value = 35229.285
$ws.Range("J106").Value = 35229.285
$ws.Range("L106").Value = 35229.285
$ws.Range("N106").Value = -37753.285
$ws.Range("H136").Value = 5190.24
$ws.Range("I136").Value = 4024.8333
$ws.Range("J136").Value = 6266
$ws.Range("K136").Value = 12074.4999
$ws.Range("L136").Value = 18798
$ws.Range("M136").Value = -9524.499899999999
$ws.Range("N136").Value = -23898

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 4763403.5
$ws.Range("I100").Value = 7938395.5
$ws.Range("J100").Value = 915
$ws.Range("K100").Value = 15876791
$ws.Range("L100").Value = 1830
$ws.Range("M100").Value = -15876250
$ws.Range("N100").Value = -2912
$ws.Range("H122").Value = 3061.1765
$ws.Range("I122").Value = 2668.5386
$ws.Range("J122").Value = 4337.25
$ws.Range("K122").Value = 8005.6158
$ws.Range("L122").Value = 13011.75
$ws.Range("M122").Value = -5555.6158
$ws.Range("N122").Value = -17911.75
$ws.Range("H126").Value = 1054.4667
$ws.Range("I126").Value = 1054.4667
$ws.Range("K126").Value = 3163.4001
$ws.Range("M126").Value = -693.4000999999998
$ws.Range("H132").Value = 1176063.1
$ws.Range("I132").Value = 814.74194
$ws.Range("J132").Value = 7248180
$ws.Range("K132").Value = 2444.22582
$ws.Range("L132").Value = 21744540
$ws.Range("M132").Value = 85.77417999999989
$ws.Range("N132").Value = -21749600
